# Updated cryptos list with GitHub Actions - refresh prices/volumes and
# re-sort a few coins (Dai/BitcoinCash, HuobiToken/LidoDAOToken,
# MXToken/FraxShare swapped order) and replace EnergySwap with BabyDogeCoin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that look like plain numbers ("210.67", "1.00", ...)
# are written with a leading apostrophe so Excel keeps them as text
# (matching the workbook's original inlineStr string cells) instead of
# silently converting them to numeric values.

# Row 2
$ws.Range("D2").Value = "26.498.03"
$ws.Range("E2").Value = "  -2.53%  "

# Row 3
$ws.Range("D3").Value = "1.581.66"
$ws.Range("E3").Value = "  -3.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").Value = "'210.67"
$ws.Range("E5").Value = "  -2.70%  "

# Row 6
$ws.Range("E6").Value = "  -2.13%  "

# Row 7
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("E8").Value = "  -2.30%  "

# Row 9
$ws.Range("E9").Value = "  -1.18%  "

# Row 10
$ws.Range("E10").Value = "  -4.02%  "

# Row 11
$ws.Range("E11").Value = "  -2.17%  "

# Row 12
$ws.Range("D12").Value = "1.802.05"
$ws.Range("E12").Value = "  -3.18%  "

# Row 13
$ws.Range("D13").Value = "1.587.63"
$ws.Range("E13").Value = "  -2.86%  "

# Row 14
$ws.Range("E14").Value = "  -1.67%  "

# Row 15
$ws.Range("E15").Value = "  -2.83%  "

# Row 16
$ws.Range("D16").Value = "'63.80"
$ws.Range("E16").Value = "  -1.50%  "

# Row 17
$ws.Range("D17").Value = "26.534.59"
$ws.Range("E17").Value = "  -2.26%  "

# Row 18
$ws.Range("E18").Value = "  -0.82%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'208.17"
$ws.Range("E19").Value = "  -3.13%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  +0.30%  "

# Row 21
$ws.Range("E21").Value = "  -3.12%  "

# Row 22
$ws.Range("D22").Value = "'4.24"
$ws.Range("E22").Value = "  -3.58%  "

# Row 23
$ws.Range("E23").Value = "  -4.90%  "

# Row 24
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  -2.04%  "

# Row 25
$ws.Range("D25").Value = "'146.23"
$ws.Range("E25").Value = "  -1.37%  "

# Row 26
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("E27").Value = "  +1.76%  "

# Row 28
$ws.Range("E28").Value = "  -4.44%  "

# Row 29
$ws.Range("D29").Value = "'15.24"
$ws.Range("E29").Value = "  -2.08%  "

# Row 30
$ws.Range("D30").Value = "'0.0500"
$ws.Range("E30").Value = "  -0.83%  "

# Row 31
$ws.Range("E31").Value = "  -2.27%  "

# Row 32
$ws.Range("E32").Value = "  -3.83%  "

# Row 33
$ws.Range("E33").Value = "  +21.97%  "

# Row 34
$ws.Range("E34").Value = "  -2.51%  "

# Row 35
$ws.Range("D35").Value = "1.305.84"
$ws.Range("E35").Value = "  -0.97%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = "  -0.58%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.50"
$ws.Range("E37").Value = "  -3.80%  "

# Row 38
$ws.Range("E38").Value = "  -0.99%  "

# Row 39
$ws.Range("E39").Value = "  -3.19%  "

# Row 40
$ws.Range("E40").Value = "  +0.27%  "

# Row 41
$ws.Range("D41").Value = "'0.782"
$ws.Range("E41").Value = "  -2.72%  "

# Row 42
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  -4.33%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.28"
$ws.Range("E43").Value = "  +0.87%  "

# Row 44
$ws.Range("D44").Value = "'62.65"
$ws.Range("E44").Value = "  -1.86%  "

# Row 45
$ws.Range("D45").Value = "1.716.02"
$ws.Range("E45").Value = "  -2.95%  "

# Row 46
$ws.Range("D46").Value = "'88.84"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47
$ws.Range("E47").Value = "  +0.35%  "

# Row 48
$ws.Range("D48").Value = "'0.831"
$ws.Range("E48").Value = "  +5.29%  "

# Row 49
$ws.Range("E49").Value = "  -1.80%  "

# Row 50
$ws.Range("D50").Value = "'0.0979"
$ws.Range("E50").Value = "  +3.34%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0961"
$ws.Range("E51").Value = "  -8.94%  "
